# Update the example register name on the "aggressive" sheet (A2) from the
# old "soc.north.vpupll" placeholder to the new "cdie.atom0.pma_gpsb" one.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("aggressive")
$ws.Range("A2").Value = "cdie.atom0.pma_gpsb"

# Widen column A slightly (stored OOXML width 14 == ColumnWidth 14 - 5/6).
$ws.Columns.Item(1).ColumnWidth = 13.166666666666666

# Move the sheet's active selection from F6 to D8.
$ws.Range("D8").Select() | Out-Null
